# Remove the trailing "Ver no Jupiter ..." and "© 2020 ... Jekyll ..."
# paragraphs (plus the blank paragraph that separates them from the
# bibliography entry above), leaving the bibliography's last entry
# followed directly by the document's closing blank paragraph and the
# page-break paragraph.

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
$jupiterIndex = -1
$jekyllIndex = -1

for ($i = 1; $i -le $count; $i++) {
    $text = $d.Paragraphs($i).Range.Text
    if ($jupiterIndex -eq -1 -and $text -like "*Ver no Jupiter*") {
        $jupiterIndex = $i
    }
    if ($text -like "*Powered by Jekyll*") {
        $jekyllIndex = $i
    }
}

if ($jupiterIndex -ne -1 -and $jekyllIndex -ne -1) {
    # Start just before the "Ver no Jupiter" paragraph so the blank
    # paragraph preceding it is removed too; end at the end of the
    # "Powered by Jekyll" (copyright) paragraph, consuming its mark.
    $startPara = $d.Paragraphs($jupiterIndex - 1)
    $endPara = $d.Paragraphs($jekyllIndex)

    $range = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $range.Delete()
}
